$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)

$ws.Range("D2").Value = '29.062.83'
$ws.Range("E2").Value = '  +0.06%  '

$ws.Range("D3").Value = '1.834.52'
$ws.Range("E3").Value = '  +0.29%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.57'
$ws.Range("E5").Value = '  +1.64%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6335'
$ws.Range("E6").Value = '  +2.07%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07546'
$ws.Range("E8").Value = '  +0.75%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2947'
$ws.Range("E9").Value = '  +1.21%  '

$ws.Range("E10").Value = '  +1.18%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07742'
$ws.Range("E11").Value = '  +1.55%  '

$ws.Range("D12").Value = '1.839.08'
$ws.Range("E12").Value = '  +0.67%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.007'
$ws.Range("E13").Value = '  +1.14%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6708'
$ws.Range("E14").Value = '  +1.22%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '83.27'
$ws.Range("E15").Value = '  +1.46%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009634'
$ws.Range("E16").Value = '  +5.42%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.090'
$ws.Range("E17").Value = '  +1.83%  '

$ws.Range("D18").Value = '29.098.86'
$ws.Range("E18").Value = '  +0.20%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.59'
$ws.Range("E19").Value = '  +2.18%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '226.35'
$ws.Range("E20").Value = '  +0.75%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9997'
$ws.Range("E21").Value = '  -0.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.202'
$ws.Range("E22").Value = '  +0.37%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.000'
$ws.Range("E23").Value = '  +0.05%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '160.69'
$ws.Range("E24").Value = '  +0.79%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1402'
$ws.Range("E25").Value = '  +3.54%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.549'
$ws.Range("E26").Value = '  +1.71%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.95'
$ws.Range("E27").Value = '  +0.81%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.503'
$ws.Range("E28").Value = '  +0.64%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.127'
$ws.Range("E29").Value = '  +2.00%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.074'
$ws.Range("E30").Value = '  +1.11%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.204'

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05389'
$ws.Range("E32").Value = '  +3.36%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.866'
$ws.Range("E33").Value = '  +1.95%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7455'
$ws.Range("E34").Value = '  +1.90%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.143'
$ws.Range("E35").Value = '  -0.58%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.658'
$ws.Range("E36").Value = '  +0.55%  '

$ws.Range("D37").Value = '1.243.14'
$ws.Range("E37").Value = '  -2.26%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.757'
$ws.Range("E38").Value = '  +0.35%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01789'
$ws.Range("E39").Value = '  +0.54%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.621'
$ws.Range("E40").Value = '  +4.99%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9064'
$ws.Range("E41").Value = '  +1.46%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.000'
$ws.Range("E42").Value = '  -0.06%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.92'
$ws.Range("E43").Value = '  +0.06%  '

$ws.Range("D44").Value = '1.985.22'
$ws.Range("E44").Value = '  +0.52%  '

$ws.Range("E45").Value = '  +3.42%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.88'
$ws.Range("E46").Value = '  +2.55%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5113'
$ws.Range("E47").Value = '  -0.09%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4093'
$ws.Range("E48").Value = '  +3.41%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.054'
$ws.Range("E49").Value = '  +2.10%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.781'
$ws.Range("E50").Value = '  +2.11%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05776'
$ws.Range("E51").Value = '  +0.42%  '
